$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 274
$ws1.Range("F5").Value = 7627
$ws1.Range("F6").Value = 5547
$ws1.Range("F8").Value = 70
$ws1.Range("F9").Value = 10
$ws1.Range("F11").Value = 245
$ws1.Range("F12").Value = 182

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 274
$ws4.Range("F5").Value = 7627
$ws4.Range("F6").Value = 5547
$ws4.Range("F8").Value = 70
$ws4.Range("F9").Value = 10
$ws4.Range("F11").Value = 245
$ws4.Range("F14").Value = 182
